$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 75064.19
$ws.Range("J17").Value = 75064.19
$ws.Range("L17").Value = 225192.57
$ws.Range("N17").Value = -225528.57
$ws.Range("H74").Value = 4967.25
$ws.Range("I74").Value = 4967.25
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4967.25
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4031.25
$ws.Range("N74").Value = $null
$ws.Range("H76").Value = 55616730
$ws.Range("I76").Value = 129793.125
$ws.Range("K76").Value = 129793.125
$ws.Range("M76").Value = -129478.125
$ws.Range("H77").Value = 4967.25
$ws.Range("I77").Value = 4967.25
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 24836.25
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -20156.25
$ws.Range("N77").Value = $null
$ws.Range("H79").Value = 55616730
$ws.Range("I79").Value = 129793.125
$ws.Range("K79").Value = 129793.125
$ws.Range("M79").Value = -128701.125
$ws.Range("H99").Value = 84489.664
$ws.Range("I99").Value = 443.4
$ws.Range("J99").Value = 144522.72
$ws.Range("K99").Value = 1330.2
$ws.Range("L99").Value = 433568.16
$ws.Range("M99").Value = 167.8000000000002
$ws.Range("N99").Value = -436564.16
$ws.Range("H100").Value = 8900
$ws.Range("I100").Value = 8900
$ws.Range("K100").Value = 8900
$ws.Range("M100").Value = -8359
$ws.Range("H137").Value = 13185.806
$ws.Range("I137").Value = 1360.2727
$ws.Range("J137").Value = 18389.04
$ws.Range("K137").Value = 4080.8181
$ws.Range("L137").Value = 55167.12
$ws.Range("M137").Value = -1530.8181
$ws.Range("N137").Value = -60267.12
$ws.Range("H138").Value = 5157.9824
$ws.Range("J138").Value = 6241.171
$ws.Range("L138").Value = 18723.513
$ws.Range("N138").Value = -29003.513

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1589.8644
$ws.Range("I32").Value = 1685.6923
$ws.Range("J32").Value = 878
$ws.Range("K32").Value = 1685.6923
$ws.Range("L32").Value = 878
$ws.Range("M32").Value = -1398.6923
$ws.Range("N32").Value = -1452
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null
$ws.Range("H63").Value = 7327.231
$ws.Range("I63").Value = 4950
$ws.Range("K63").Value = 4950
$ws.Range("M63").Value = -4264
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null
$ws.Range("H66").Value = 7327.231
$ws.Range("I66").Value = 4950
$ws.Range("K66").Value = 24750
$ws.Range("M66").Value = -21318

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 657.93335
$ws.Range("I94").Value = 577.0909
$ws.Range("J94").Value = 880.25
$ws.Range("K94").Value = 577.0909
$ws.Range("L94").Value = 880.25
$ws.Range("M94").Value = -126.0909
$ws.Range("N94").Value = -1782.25
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").Value = $null
$ws.Range("H132").Value = 49999.332
$ws.Range("J132").Value = 49999.332
$ws.Range("L132").Value = 49999.332
$ws.Range("N132").Value = -60119.332

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 355822.66
$ws.Range("I31").Value = 771367.25
$ws.Range("J31").Value = 17230.777
$ws.Range("K31").Value = 771367.25
$ws.Range("L31").Value = 17230.777
$ws.Range("M31").Value = -771072.25
$ws.Range("N31").Value = -17820.777
$ws.Range("H34").Value = 355822.66
$ws.Range("I34").Value = 771367.25
$ws.Range("J34").Value = 17230.777
$ws.Range("K34").Value = 771367.25
$ws.Range("L34").Value = 17230.777
$ws.Range("M34").Value = -771165.25
$ws.Range("N34").Value = -17634.777
$ws.Range("H68").Value = 85600
$ws.Range("J68").Value = 85800
$ws.Range("L68").Value = 85800
$ws.Range("N68").Value = -87298
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H71").Value = 85600
$ws.Range("J71").Value = 85800
$ws.Range("L71").Value = 257400
$ws.Range("N71").Value = -264888
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H74").Value = 60000
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61748
$ws.Range("H77").Value = 60000
$ws.Range("J77").Value = 60000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -188736
$ws.Range("H92").Value = 49884.5
$ws.Range("J92").Value = 49884.5
$ws.Range("L92").Value = 49884.5
$ws.Range("N92").Value = -54876.5
$ws.Range("H99").Value = 8500
$ws.Range("J99").Value = 8500
$ws.Range("L99").Value = 8500
$ws.Range("N99").Value = -11496
$ws.Range("H105").Value = 2198.25
$ws.Range("J105").Value = 4995
$ws.Range("L105").Value = 4995
$ws.Range("N105").Value = -8489
$ws.Range("H122").Value = 3899.6667
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3899.6667
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11699.0001
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = -16599.0001
$ws.Range("H126").Value = 8500
$ws.Range("J126").Value = 8500
$ws.Range("L126").Value = 25500
$ws.Range("N126").Value = -30440

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 18500
$ws.Range("J101").Value = 18500
$ws.Range("L101").Value = 55500
$ws.Range("N101").Value = -60368
$ws.Range("H114").Value = 367.22223
$ws.Range("I114").Value = 392.66666
$ws.Range("K114").Value = 1177.99998
$ws.Range("M114").Value = 2076.00002
$ws.Range("H117").Value = 1085
$ws.Range("J117").Value = 1500
$ws.Range("L117").Value = 4500
$ws.Range("N117").Value = -11384

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 1104.8889
$ws.Range("I70").Value = 1102.6
$ws.Range("K70").Value = 1102.6
$ws.Range("M70").Value = -832.5999999999999
$ws.Range("H73").Value = 1104.8889
$ws.Range("I73").Value = 1102.6
$ws.Range("K73").Value = 1102.6
$ws.Range("M73").Value = -166.5999999999999
$ws.Range("H122").Value = 309897.5
$ws.Range("I122").Value = 383246.62
$ws.Range("K122").Value = 1149739.86
$ws.Range("M122").Value = -1147289.86

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1688168.9
$ws.Range("I7").Value = 3339669.2
$ws.Range("J7").Value = 36668.332
$ws.Range("K7").Value = 3339669.2
$ws.Range("L7").Value = 36668.332
$ws.Range("M7").Value = -3339557.2
$ws.Range("N7").Value = -36892.332
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H100").Value = 146562.38
$ws.Range("I100").Value = 363333
$ws.Range("K100").Value = 363333
$ws.Range("M100").Value = -362792
$ws.Range("H126").Value = 1688168.9
$ws.Range("I126").Value = 3339669.2
$ws.Range("J126").Value = 36668.332
$ws.Range("K126").Value = 10019007.6
$ws.Range("L126").Value = 110004.996
$ws.Range("M126").Value = -10016537.6
$ws.Range("N126").Value = -114944.996
$ws.Range("H132").Value = 3359.7754
$ws.Range("I132").Value = 2706.3142
$ws.Range("J132").Value = 4993.4287
$ws.Range("K132").Value = 8118.942599999999
$ws.Range("L132").Value = 14980.2861
$ws.Range("M132").Value = -5588.942599999999
$ws.Range("N132").Value = -20040.2861
$ws.Range("H133").Value = 56666.445
$ws.Range("J133").Value = 56666.445
$ws.Range("L133").Value = 56666.445
$ws.Range("N133").Value = -61726.445

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = $null
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = $null
$ws.Range("H98").Value = 999999
$ws.Range("J98").Value = 999999
$ws.Range("L98").Value = 999999
$ws.Range("N98").Value = -1005989
$ws.Range("H107").Value = 46215.453
$ws.Range("I107").Value = 63124.188
$ws.Range("K107").Value = 189372.564
$ws.Range("M107").Value = -187452.564
$ws.Range("H122").Value = 3704.8438
$ws.Range("I122").Value = 3253.8845
$ws.Range("K122").Value = 9761.6535
$ws.Range("M122").Value = -7311.6535
